# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) timestamps on the
# per-language report rows (row 2 of each language sheet).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-06 04:01:03"
$wsZhCn.Range("G2").Value = "2016-02-06 04:01:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-06 04:01:14"
$wsDeDe.Range("G2").Value = "2016-02-06 04:02:05"
